$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 3.8 e享课堂项目预算: update programmer cost figures (row 12-14, column F "总金额")
$ws.Range("F12").Value = "￥13000"
$ws.Range("F13").Value = "￥10000"
$ws.Range("F14").Value = "￥10000"

# Move the active selection to match the saved workbook state
$ws.Range("E16").Select() | Out-Null
